$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.319.63"
$ws.Range("E2").Value = "  +5.32%  "
$ws.Range("D3").Value = "3.162.58"
$ws.Range("E3").Value = "  +2.96%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'398.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.18%  "
$ws.Range("D6").Value = "'109.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.25%  "
$ws.Range("D7").Value = "'0.549"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.98%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +5.06%  "
$ws.Range("D10").Value = "'38.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.63%  "
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("E12").Value = "  +1.72%  "
$ws.Range("D13").Value = "3.660.21"
$ws.Range("E13").Value = "  +3.03%  "
$ws.Range("D14").Value = "'19.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.32%  "
$ws.Range("D15").Value = "'8.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.05%  "
$ws.Range("E16").Value = "  +8.16%  "
$ws.Range("D17").Value = "3.164.82"
$ws.Range("E17").Value = "  +3.40%  "
$ws.Range("D18").Value = "'10.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.36%  "
$ws.Range("D19").Value = "54.044.89"
$ws.Range("E19").Value = "  +4.72%  "
$ws.Range("E20").Value = "  +3.95%  "
$ws.Range("E21").Value = "  +2.81%  "
$ws.Range("D22").Value = "0.0₃0979"
$ws.Range("E22").Value = "  +1.03%  "
$ws.Range("D23").Value = "'71.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.21%  "
$ws.Range("D24").Value = "'271.68"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.16%  "
$ws.Range("D25").Value = "'3.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.53%  "
$ws.Range("D26").Value = "'8.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.22%  "
$ws.Range("D27").Value = "'27.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.57%  "
$ws.Range("D28").Value = "'7.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.83%  "
$ws.Range("D29").Value = "'0.170"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  +3.58%  "
$ws.Range("D32").Value = "'10.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.89%  "
$ws.Range("D33").Value = "'0.0506"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +13.07%  "
$ws.Range("D34").Value = "'37.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.30%  "
$ws.Range("E35").Value = "  +0.40%  "
$ws.Range("D36").Value = "'50.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("D37").Value = "'3.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.15%  "
$ws.Range("D38").Value = "'0.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("D39").Value = "'2.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.66%  "
$ws.Range("E40").Value = "  +9.29%  "
$ws.Range("D41").Value = "'0.292"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.02%  "
$ws.Range("D42").Value = "'17.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.46%  "
$ws.Range("E43").Value = "  +1.15%  "
$ws.Range("D44").Value = "'130.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.16%  "
$ws.Range("D45").Value = "'0.117"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.23%  "
$ws.Range("D46").Value = "'22.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.57%  "
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("D48").Value = "'2.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("D49").Value = "2.082.57"
$ws.Range("E49").Value = "  +2.30%  "
$ws.Range("D50").Value = "'0.0344"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.78%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "'5.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.16%  "
